$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new rows of test data (testT4239) below the existing rows
$ws.Range("A16").Value = "testT4239"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "YES"

$ws.Range("A17").Value = "testT4239"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "YES"

# Update the active selection to match the newly added row
$ws.Range("A17:C17").Select()
